$d = $word.ActiveDocument

$replacements = @(
    @("521×3=", "393×9="),
    @("730×3=", "229×2="),
    @("812×6=", "988×8="),
    @("143×5=", "902×7="),
    @("954×5=", "349×2="),
    @("191×2=", "912×6="),
    @("736×2=", "946×2="),
    @("705×7=", "164×9="),
    @("337×4=", "237×8="),
    @("821×2=", "641×6="),
    @("451×8=", "992×2="),
    @("228×9=", "660×2="),
    @("829×8=", "777×7="),
    @("229×3=", "859×2="),
    @("798×8=", "591×7="),
    @("497×4=", "657×2="),
    @("317×8=", "772×6="),
    @("441×3=", "154×2="),
    @("725×3=", "937×4="),
    @("382×8=", "362×2="),
    @("284×7=", "386×4="),
    @("233×7=", "377×6="),
    @("959×8=", "322×8="),
    @("443×9=", "421×3="),
    @("869×4=", "222×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
